# ------------------------------------------------------------------
# Applies the "Improve pairing generator, Flight constructor and other
# small details" commit to the workbook:
#   - Flights sheet: new/extended flight pairing data (rows 2-20),
#     wider C/D date columns with a custom date/time number format,
#     a new column E width tweak, and refreshed selection/active sheet.
#   - Parameters sheet loses the "tabSelected" flag (Flights becomes
#     the selected tab instead).
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$flights = $wb.Worksheets.Item("Flights")
$params  = $wb.Worksheets.Item("Parameters")

# --- New flight pairing data (origin, destination, start, end) -----
$flights.Cells.Item(2,1).Value = "DFW"
$flights.Cells.Item(2,2).Value = "JFK"
$flights.Cells.Item(2,3).Value = 45451.25
$flights.Cells.Item(2,4).Value = 45451.375

$flights.Cells.Item(3,1).Value = "JFK"
$flights.Cells.Item(3,2).Value = "LAX"
$flights.Cells.Item(3,3).Value = 45451.395833333336
$flights.Cells.Item(3,4).Value = 45451.520833333336

$flights.Cells.Item(4,1).Value = "LAX"
$flights.Cells.Item(4,2).Value = "ORD"
$flights.Cells.Item(4,3).Value = 45451.541666666664
$flights.Cells.Item(4,4).Value = 45451.708333333336

$flights.Cells.Item(5,1).Value = "ORD"
$flights.Cells.Item(5,2).Value = "DFW"
$flights.Cells.Item(5,3).Value = 45451.729166666664
$flights.Cells.Item(5,4).Value = 45451.833333333336

$flights.Cells.Item(6,1).Value = "DFW"
$flights.Cells.Item(6,2).Value = "LAX"
$flights.Cells.Item(6,3).Value = 45452.25
$flights.Cells.Item(6,4).Value = 45452.375

$flights.Cells.Item(7,1).Value = "LAX"
$flights.Cells.Item(7,2).Value = "JFK"
$flights.Cells.Item(7,3).Value = 45452.395833333336
$flights.Cells.Item(7,4).Value = 45452.625

$flights.Cells.Item(8,1).Value = "JFK"
$flights.Cells.Item(8,2).Value = "ORD"
$flights.Cells.Item(8,3).Value = 45452.645833333336
$flights.Cells.Item(8,4).Value = 45452.75

$flights.Cells.Item(9,1).Value = "ORD"
$flights.Cells.Item(9,2).Value = "LAX"
$flights.Cells.Item(9,3).Value = 45452.770833333336
$flights.Cells.Item(9,4).Value = 45452.895833333336

$flights.Cells.Item(10,1).Value = "LAX"
$flights.Cells.Item(10,2).Value = "DFW"
$flights.Cells.Item(10,3).Value = 45452.916666666664
$flights.Cells.Item(10,4).Value = 45453.041666666664

$flights.Cells.Item(11,1).Value = "DFW"
$flights.Cells.Item(11,2).Value = "MIA"
$flights.Cells.Item(11,3).Value = 45453.25
$flights.Cells.Item(11,4).Value = 45453.354166666664

$flights.Cells.Item(12,1).Value = "MIA"
$flights.Cells.Item(12,2).Value = "ATL"
$flights.Cells.Item(12,3).Value = 45453.375
$flights.Cells.Item(12,4).Value = 45453.458333333336

$flights.Cells.Item(13,1).Value = "ATL"
$flights.Cells.Item(13,2).Value = "ORD"
$flights.Cells.Item(13,3).Value = 45453.479166666664
$flights.Cells.Item(13,4).Value = 45453.5625

$flights.Cells.Item(14,1).Value = "ORD"
$flights.Cells.Item(14,2).Value = "SFO"
$flights.Cells.Item(14,3).Value = 45453.583333333336
$flights.Cells.Item(14,4).Value = 45453.708333333336

$flights.Cells.Item(15,1).Value = "SFO"
$flights.Cells.Item(15,2).Value = "LAX"
$flights.Cells.Item(15,3).Value = 45453.729166666664
$flights.Cells.Item(15,4).Value = 45453.770833333336

$flights.Cells.Item(16,1).Value = "LAX"
$flights.Cells.Item(16,2).Value = "SEA"
$flights.Cells.Item(16,3).Value = 45453.791666666664
$flights.Cells.Item(16,4).Value = 45453.895833333336

$flights.Cells.Item(17,1).Value = "SEA"
$flights.Cells.Item(17,2).Value = "ORD"
$flights.Cells.Item(17,3).Value = 45454.25
$flights.Cells.Item(17,4).Value = 45454.458333333336

$flights.Cells.Item(18,1).Value = "ORD"
$flights.Cells.Item(18,2).Value = "JFK"
$flights.Cells.Item(18,3).Value = 45454.479166666664
$flights.Cells.Item(18,4).Value = 45454.5625

$flights.Cells.Item(19,1).Value = "JFK"
$flights.Cells.Item(19,2).Value = "BOS"
$flights.Cells.Item(19,3).Value = 45454.583333333336
$flights.Cells.Item(19,4).Value = 45454.645833333336

$flights.Cells.Item(20,1).Value = "BOS"
$flights.Cells.Item(20,2).Value = "DFW"
$flights.Cells.Item(20,3).Value = 45454.666666666664
$flights.Cells.Item(20,4).Value = 45454.8125

# --- Formatting: start/end columns get a custom date/time format ---
$flights.Range("C1:D20").NumberFormat = "yyyy/mm/dd\ hh:mm"

# --- Column widths (C, D widened; new explicit width for E) --------
$flights.Columns.Item(3).ColumnWidth = 21.25
$flights.Columns.Item(4).ColumnWidth = 19.583333333333332
$flights.Columns.Item(5).ColumnWidth = 12.916666666666666

# --- Selection / active sheet tweaks --------------------------------
$flights.Range("E9").Select()
$flights.Activate()

Write-Host "done"
